$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lealtades")

# Update the data-binding placeholder in F4 from {{item.IdListaPrecios}} to {{item.PrecioLista}}
$ws.Range("F4").Value = "{{item.PrecioLista}}"

# Update the active selection to match the new cursor position (F4)
$ws.Activate()
$ws.Range("F4").Select()
